$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match header style to the existing header cells (A1:E1)
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Rows where the KNN_Outliers_MAD flag is TRUE
$trueRows = @(12, 15, 19)

for ($r = 2; $r -le 21; $r++) {
    if ($trueRows -contains $r) {
        $ws.Cells.Item($r, 6).Value = $true
    } else {
        $ws.Cells.Item($r, 6).Value = $false
    }
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
